$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-01-31 Friday" "2025-02-01 Saturday"
Replace-Text "693÷3=" "734÷5="
Replace-Text "964÷7=" "352÷3="
Replace-Text "118÷9=" "275÷4="
Replace-Text "581÷5=" "946÷8="
Replace-Text "912÷2=" "216÷4="
Replace-Text "258÷4=" "917÷6="
Replace-Text "682÷5=" "912÷5="
Replace-Text "382÷3=" "719÷4="
Replace-Text "850÷7=" "659÷5="
Replace-Text "808÷6=" "366÷2="
Replace-Text "394÷7=" "350÷2="
Replace-Text "807÷3=" "857÷8="
Replace-Text "671÷3=" "247÷3="
Replace-Text "720÷2=" "909÷4="
Replace-Text "487÷9=" "443÷9="
Replace-Text "853÷6=" "250÷4="
Replace-Text "545÷7=" "437÷7="
Replace-Text "933÷2=" "106÷4="
Replace-Text "337÷2=" "376÷3="
Replace-Text "991÷5=" "797÷9="
Replace-Text "408÷5=" "430÷8="
Replace-Text "647÷8=" "266÷9="
Replace-Text "116÷2=" "650÷4="
Replace-Text "925÷3=" "495÷6="
Replace-Text "674÷9=" "678÷2="
